$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns I and J, with the same bold/centered/bordered style
# used by the existing headers (e.g. H1 "IP"). Copy H1's formatting over
# before writing the new header text so the new cells pick up style index 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2 (plain, unstyled numeric cells like C2:H2).
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7
